$wb = $excel.ActiveWorkbook

# weibull
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -1.91600227134157
$ws.Range("C2").Value = 0.257928758189887
$ws.Range("B3").Value = -0.139560850406393
$ws.Range("C3").Value = 0.105925170492656

# lognormal
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 1.02869503297525
$ws.Range("C2").Value = 0.265769224642769
$ws.Range("B3").Value = -0.717688599529452
$ws.Range("C3").Value = 0.103533830249769

# llogis
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -1.46716693033476
$ws.Range("C2").Value = 0.287705749115125
$ws.Range("B3").Value = 1.26059972867232
$ws.Range("C3").Value = 0.170502953872554

# gompertz
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -1.77351822022407
$ws.Range("C2").Value = 0.320554106187069
$ws.Range("B3").Value = -0.0531899742245637
$ws.Range("C3").Value = 0.0449845807907303

# weibull cov
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0665272443013772
$ws.Range("B2").Value = -0.0112609082153799
$ws.Range("A3").Value = -0.0112609082153799
$ws.Range("B3").Value = 0.0112201417438982

# lognormal cov
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0706332807672187
$ws.Range("B2").Value = -0.0171998828454053
$ws.Range("A3").Value = -0.0171998828454053
$ws.Range("B3").Value = 0.010719254006188

# llogis cov
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0827745980738954
$ws.Range("B2").Value = 0.0044798074636298
$ws.Range("A3").Value = 0.0044798074636298
$ws.Range("B3").Value = 0.0290712572792663

# gompertz cov
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.10275493499339
$ws.Range("B2").Value = -0.0106890316796556
$ws.Range("A3").Value = -0.0106890316796556
$ws.Range("B3").Value = 0.00202361250891774
